$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new cell values in the exact order needed so the shared-string
# table grows in the same sequence as the target workbook ---
$ws.Range("B32").Value = "T15"
$ws.Range("K32").Value = "Hierarchic workflow test."
$ws.Range("A32").Value = "015-shell-hierarchic-workflow"
$ws.Range("A31").Value = "014-planner-performance"
$ws.Range("A30").Value = "013-pegasus-mpi-cluster"
$ws.Range("B31").Value = "T014"
$ws.Range("B30").Value = "TESTMPIDAG"
$ws.Range("C32").Value = "Shell"

# --- Fix up formatting on column A / B cells that need a non-default style
# (copying pulls in the template cell's style; re-setting the value afterwards
# keeps that style because the cell/row already exists) ---
$ws.Range("A6").Copy($ws.Range("A30"))
$ws.Range("A30").Value = "013-pegasus-mpi-cluster"

$ws.Range("A6").Copy($ws.Range("A31"))
$ws.Range("A31").Value = "014-planner-performance"

$ws.Range("A6").Copy($ws.Range("A32"))
$ws.Range("A32").Value = "015-shell-hierarchic-workflow"

$ws.Range("B6").Copy($ws.Range("B30"))
$ws.Range("B30").Value = "TESTMPIDAG"

$ws.Range("B6").Copy($ws.Range("B32"))
$ws.Range("B32").Value = "T15"

# --- Column B gets a bit wider to fit the new "TESTMPIDAG" entry ---
$ws.Columns("B").ColumnWidth = 11.5

# --- Match the saved selection state ---
[void]$ws.Range("C8").Select()
